# Update "想去人数" (want-to-go count) column F across the four sheets to
# reflect the refreshed scrape output (gh-pages data refresh @ 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) -------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 12323
$ws.Range("F3").Value = 6920
$ws.Range("F5").Value = 21
$ws.Range("F7").Value = 263
$ws.Range("F11").Value = 121
$ws.Range("F12").Value = 322
$ws.Range("F13").Value = 973
$ws.Range("F14").Value = 3701
$ws.Range("F16").Value = 997
$ws.Range("F18").Value = 215
$ws.Range("F19").Value = 342
$ws.Range("F23").Value = 74
$ws.Range("F25").Value = 5124
$ws.Range("F26").Value = 57
$ws.Range("F27").Value = 1338
$ws.Range("F28").Value = 267
$ws.Range("F29").Value = 807
$ws.Range("F30").Value = 1284
$ws.Range("F31").Value = 572

# --- Sheet "演出" (Performance) -------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 3725
$ws.Range("F7").Value = 27

# --- Sheet "本地生活" (Local life) -----------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9176
$ws.Range("F3").Value = 537
$ws.Range("F4").Value = 1916

# --- Sheet "全部类型" (All types) -------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9176
$ws.Range("F3").Value = 537
$ws.Range("F4").Value = 1916
$ws.Range("F5").Value = 12323
$ws.Range("F6").Value = 6920
$ws.Range("F8").Value = 3725
$ws.Range("F10").Value = 21
$ws.Range("F12").Value = 263
$ws.Range("F16").Value = 121
$ws.Range("F17").Value = 322
$ws.Range("F18").Value = 973
$ws.Range("F19").Value = 3701
$ws.Range("F21").Value = 997
$ws.Range("F23").Value = 215
$ws.Range("F24").Value = 342
$ws.Range("F33").Value = 5124
$ws.Range("F34").Value = 57
$ws.Range("F35").Value = 1338
$ws.Range("F38").Value = 267
$ws.Range("F40").Value = 807
$ws.Range("F41").Value = 1284
$ws.Range("F42").Value = 572
